$wb = $excel.ActiveWorkbook

# Existing sheet that the new sheet will be inserted before
$ws19T4 = $wb.Worksheets.Item("19T4")

# Insert the new sheet "18R1" right before "19T4" (so order becomes 19S1, 22T1, 18R1, 19T4)
$newSheet = $wb.Worksheets.Add($ws19T4)
$newSheet.Name = "18R1"

# Header row (row 1): bold style matching existing header style used across the workbook
$newSheet.Range("A1:I1").Font.Color = 0
$newSheet.Range("A1:I1").Font.Bold = $true

# Data rows (rows 2-5): plain style matching existing data style used across the workbook
$newSheet.Range("A2:I5").Font.Color = 0

# Header values
$newSheet.Range("A1").Value = "Analysis Population"
$newSheet.Range("B1").Value = "Valid Tests Size"
$newSheet.Range("C1").Value = "Positive Tests"
$newSheet.Range("D1").Value = "Median Percent (%)"

# Data values
$newSheet.Range("A2").Value = "Indiviuals older than 15y"
$newSheet.Range("B2").Value = 1319
$newSheet.Range("C2").Value = 286
$newSheet.Range("D2").Value = 36.8

# Selection / active cell for the new sheet view
$newSheet.Range("F19:F20").Select()
